# Actualización automática 2025-11-12 17:30:08
# Updates the PORCELANATO sales figure for client "GANCHOZO CEDEÑO YURI MERCEDES"
# (advisor RIOS CARRION ANGEL BENIGNO) and propagates the change through the
# dependent monthly / compliance totals on the other two sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO column (M), row 14 ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M14").Value = 3080.12

# --- Sheet "VENTA MENSUAL": noviembre column (F), row 14 (client row) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F14").Value = 3080.12

# --- Sheet "VENTA MENSUAL": noviembre column (F), row 26 (TOTAL row) ---
$wsMensual.Range("F26").Value = 5495.23

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO group row (12) ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D12").Value = 5051.79
$wsCumpl.Range("E12").Value = 29649.21
$wsCumpl.Range("F12").Value = 0.1455805308204375

# --- Sheet "CUMPLIMIENTO MENSUAL": TOTAL row (14) ---
$wsCumpl.Range("D14").Value = 5495.23
$wsCumpl.Range("E14").Value = 35282.51058948192
$wsCumpl.Range("F14").Value = 0.1347605316175223
